$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.949.38"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "2.603.62"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'523.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.55%  "
$ws.Range("D6").Value = "'154.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'0.593"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.50%  "
$ws.Range("D9").Value = "'6.69"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.96%  "
$ws.Range("D10").Value = "'0.105"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").Value = "3.061.96"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").Value = "60.991.39"
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").Value = "'21.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "'0.0000141"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").Value = "2.606.42"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "'354.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.86%  "
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("E21").Value = "  +2.29%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'60.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.05%  "
$ws.Range("E24").Value = "  +2.13%  "
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").Value = "2.717.13"
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").Value = "0.0₃0848"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("D29").Value = "'7.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "'6.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.47%  "
$ws.Range("D32").Value = "'19.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("E33").Value = "  +3.17%  "
$ws.Range("D34").Value = "'148.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.21%  "
$ws.Range("D35").Value = "'4.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.91%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.13%  "
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").Value = "'0.925"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.04%  "
$ws.Range("D38").Value = "'0.878"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.69%  "
$ws.Range("E39").Value = "  +2.65%  "
$ws.Range("E40").Value = "  +2.09%  "
$ws.Range("D41").Value = "'36.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.17%  "
$ws.Range("D42").Value = "'290.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("E43").Value = "  +2.48%  "
$ws.Range("D44").Value = "'0.624"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("E45").Value = "  +1.23%  "
$ws.Range("D46").Value = "'0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "'4.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.91%  "
$ws.Range("D48").Value = "'19.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("E49").Value = "  +2.65%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("E51").Value = "  +9.23%  "
